# Apply the edits described by the commit:
#  - Shorten the "In Sheet1 the scatter plot is presented" caption to
#    "The scatter plot is presented".
#  - Remove the leftover "Sheet2" / "fit a curve" instructional text blocks
#    that no longer apply (cells are cleared but keep their formatting).
#  - Update the worksheet's saved selection to J27.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Shorten the intro caption.
$ws.Range("E8").Value = "The scatter plot is presented"

# Clear the old "In Sheet2 the plot is shown..." paragraph (3 merged lines).
$ws.Range("E14").Value = ""
$ws.Range("E15").Value = ""
$ws.Range("E16").Value = ""

# Clear the old "To fit a curve to the points..." paragraph (4 merged lines).
$ws.Range("E22").Value = ""
$ws.Range("E23").Value = ""
$ws.Range("E24").Value = ""
$ws.Range("E25").Value = ""

# Move/save the active selection to J27, as in the edited workbook.
$ws.Activate()
$ws.Range("J27").Select()
